# Apply updated probability values to Sheet1 (team-specific transition matrix)
# These changes reflect additional simulated games and refreshed optimization stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.1863636363636364
    $ws.Range("C2").Value = 0.5363636363636364
    $ws.Range("J2").Value = 0.1090909090909091
    $ws.Range("P2").Value = 0.1409090909090909
    $ws.Range("S2").Value = 0.02727272727272727
    $ws.Range("B3").Value = 0.03333333333333333
    $ws.Range("C3").Value = 0.025
    $ws.Range("J3").Value = 0.1083333333333333
    $ws.Range("P3").Value = 0.7583333333333333
    $ws.Range("S3").Value = 0.075
    $ws.Range("J4").Value = 0.2857142857142857
    $ws.Range("P4").Value = 0.6571428571428571
    $ws.Range("S4").Value = 0.05714285714285714
    $ws.Range("B6").Value = 0.04210526315789474
    $ws.Range("D6").Value = 0.01052631578947368
    $ws.Range("F6").Value = 0.05263157894736842
    $ws.Range("J6").Value = 0.3578947368421053
    $ws.Range("O6").Value = 0.04210526315789474
    $ws.Range("Q6").Value = 0.1578947368421053
    $ws.Range("R6").Value = 0.07368421052631578
    $ws.Range("S6").Value = 0.2631578947368421
    $ws.Range("B7").Value = 0.1118881118881119
    $ws.Range("D7").Value = 0.01398601398601399
    $ws.Range("F7").Value = 0.07692307692307693
    $ws.Range("J7").Value = 0.3076923076923077
    $ws.Range("O7").Value = 0.03496503496503497
    $ws.Range("Q7").Value = 0.1188811188811189
    $ws.Range("R7").Value = 0.06293706293706294
    $ws.Range("S7").Value = 0.2727272727272727
    $ws.Range("B8").Value = 0.09669211195928754
    $ws.Range("D8").Value = 0.02290076335877863
    $ws.Range("F8").Value = 0.07633587786259542
    $ws.Range("J8").Value = 0.3129770992366412
    $ws.Range("O8").Value = 0.01272264631043257
    $ws.Range("Q8").Value = 0.1450381679389313
    $ws.Range("R8").Value = 0.07633587786259542
    $ws.Range("S8").Value = 0.2569974554707379
    $ws.Range("B9").Value = 0.05084745762711865
    $ws.Range("D9").Value = 0.005649717514124294
    $ws.Range("F9").Value = 0.07909604519774012
    $ws.Range("J9").Value = 0.2598870056497175
    $ws.Range("O9").Value = 0.02824858757062147
    $ws.Range("Q9").Value = 0.1977401129943503
    $ws.Range("R9").Value = 0.0847457627118644
    $ws.Range("S9").Value = 0.2937853107344633
    $ws.Range("B10").Value = 0.06389986824769434
    $ws.Range("D10").Value = 0.01449275362318841
    $ws.Range("E10").Value = 0.0006587615283267457
    $ws.Range("F10").Value = 0.05731225296442688
    $ws.Range("J10").Value = 0.4229249011857708
    $ws.Range("O10").Value = 0.03820816864295125
    $ws.Range("Q10").Value = 0.1343873517786561
    $ws.Range("R10").Value = 0.077733860342556
    $ws.Range("S10").Value = 0.1903820816864295
    $ws.Range("G11").Value = 0.08045977011494253
    $ws.Range("J11").Value = 0.08620689655172414
    $ws.Range("K11").Value = 0.1091954022988506
    $ws.Range("L11").Value = 0.7241379310344828
    $ws.Range("G12").Value = 0.78125
    $ws.Range("J12").Value = 0.203125
    $ws.Range("K12").Value = 0.0078125
    $ws.Range("L12").Value = 0.0078125
    $ws.Range("G13").Value = 0.5675675675675675
    $ws.Range("J13").Value = 0.4054054054054054
    $ws.Range("S13").Value = 0.02702702702702703
    $ws.Range("G14").Value = 1
    $ws.Range("F15").Value = 0.01477832512315271
    $ws.Range("H15").Value = 0.1330049261083744
    $ws.Range("I15").Value = 0.09359605911330049
    $ws.Range("J15").Value = 0.458128078817734
    $ws.Range("K15").Value = 0.03448275862068965
    $ws.Range("M15").Value = 0.004926108374384237
    $ws.Range("O15").Value = 0.04433497536945813
    $ws.Range("S15").Value = 0.2167487684729064
    $ws.Range("H16").Value = 0.1468531468531468
    $ws.Range("I16").Value = 0.09090909090909091
    $ws.Range("J16").Value = 0.5314685314685315
    $ws.Range("K16").Value = 0.08391608391608392
    $ws.Range("M16").Value = 0.04895104895104895
    $ws.Range("O16").Value = 0.06993006993006994
    $ws.Range("S16").Value = 0.02797202797202797
    $ws.Range("F17").Value = 0.01479289940828402
    $ws.Range("H17").Value = 0.1804733727810651
    $ws.Range("I17").Value = 0.1183431952662722
    $ws.Range("J17").Value = 0.4911242603550296
    $ws.Range("K17").Value = 0.05621301775147929
    $ws.Range("M17").Value = 0.02662721893491124
    $ws.Range("O17").Value = 0.07100591715976332
    $ws.Range("S17").Value = 0.04142011834319527
    $ws.Range("F18").Value = 0.03208556149732621
    $ws.Range("H18").Value = 0.1657754010695187
    $ws.Range("I18").Value = 0.106951871657754
    $ws.Range("J18").Value = 0.5240641711229946
    $ws.Range("K18").Value = 0.05882352941176471
    $ws.Range("M18").Value = 0.0053475935828877
    $ws.Range("O18").Value = 0.0374331550802139
    $ws.Range("S18").Value = 0.06951871657754011
    $ws.Range("F19").Value = 0.00909090909090909
    $ws.Range("H19").Value = 0.2318181818181818
    $ws.Range("I19").Value = 0.08181818181818182
    $ws.Range("J19").Value = 0.4151515151515152
    $ws.Range("K19").Value = 0.1106060606060606
    $ws.Range("M19").Value = 0.01666666666666667
    $ws.Range("N19").Value = 0.001515151515151515
    $ws.Range("O19").Value = 0.06363636363636363
    $ws.Range("S19").Value = 0.0696969696969697

$ws.Calculate()
